# Apply BOM updates to the Power Supply Board Rev B BOM workbook.
# Resistors: MCR03EZPJ000 / MCR03EZPJ101 / MCR03EZPJ472 / MCR03EZPFX1002 /
#            MCR03EZPFX1003 -> MCR03ERTJ000 / MCR03ERTJ101 / MCR03ERTJ472 /
#            MCR03ERTF1002 / MCR03ERTF1003 (Manufacturer Part # column, F).
# Crystal X1: ECS-160-8-36CKM (ECS / XC1552CT-ND / 16MHz 8pF) replaced with
#            TXC CORPORATION 8Z-16.000MEEQ-T (Digi-Key 887-1336-1-ND,
#            16.000 MHz 10pF), Value column cleared, unit price removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resistor manufacturer part number swaps (column F) -------------------
$ws.Range("F27").Value = "MCR03ERTJ000"
$ws.Range("F31").Value = "MCR03ERTJ101"
$ws.Range("F35").Value = "MCR03ERTJ472"
$ws.Range("F37").Value = "MCR03ERTF1002"
$ws.Range("F38").Value = "MCR03ERTF1003"

# --- Crystal X1 (row 54) replacement --------------------------------------
$ws.Range("C54").ClearContents()
$ws.Range("D54").Value = "4-SMD, No Lead (DFN, LCC)"
$ws.Range("E54").Value = "TXC CORPORATION"
$ws.Range("F54").Value = "8Z-16.000MEEQ-T"
$ws.Range("H54").Value = "887-1336-1-ND"
$ws.Range("I54").Value = "CRYSTAL 16.000 MHZ 10PF SMD"
$ws.Range("K54").ClearContents()

# --- Update the saved selection to match the edit --------------------------
$ws.Range("K54").Select()
